$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, shifting existing rows 3-8 down to 4-9
$ws.Rows(3).Insert()

# Populate the new row 3 with the new data record
$ws.Range("A3").Value = 11
$ws.Range("B3").Value = "Vega Monumental Concepción"
$ws.Range("C3").Value = "Bíobío"
$ws.Range("D3").Value = 44659
$ws.Range("D3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E3").Value = 8
$ws.Range("F3").Value = 100112052
$ws.Range("G3").Value = "Albahaca"
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 90
$ws.Range("K3").Value = 2500
$ws.Range("L3").Value = 3000
$ws.Range("M3").Value = 2722
$ws.Range("N3").Value = "`$/docena de matas"
$ws.Range("O3").Value = "Región Metropolitana"
$ws.Range("P3").Value = 454
$ws.Range("Q3").Value = 6
$ws.Range("R3").Value = "Hortaliza"

Write-Host "Row inserted and populated."
